$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# This edit restructures the summary body into a bulleted/nested list and
# adds the corresponding numbering definitions (abstractNum 991, nums 1001 &
# 1002). The most reliable way to express this with the exposed object model
# is to read the canonical OOXML for the parts that need to change
# (word/document.xml and word/numbering.xml), perform targeted text surgery
# on them, and feed just those two parts back in through InsertXML. Limiting
# the round-trip to only the parts that actually change keeps every other
# part of the package byte-identical to the original.
# ---------------------------------------------------------------------------

$full = $d.Content.XML()

function Get-PartXmlData($pkg, $partName) {
    $startTag = '<pkg:part pkg:name="' + $partName + '"'
    $si = $pkg.IndexOf($startTag)
    $dataStart = $pkg.IndexOf("<pkg:xmlData>", $si) + "<pkg:xmlData>".Length
    $dataEnd = $pkg.IndexOf("</pkg:xmlData>", $dataStart)
    return $pkg.Substring($dataStart, $dataEnd - $dataStart)
}

$docXml = Get-PartXmlData $full "/word/document.xml"
$numXml = Get-PartXmlData $full "/word/numbering.xml"

# ---------------------------------------------------------------------------
# 1. Replace the body paragraphs (between <w:body> and <w:sectPr) with the
#    new restructured content.
# ---------------------------------------------------------------------------

$newBodyParagraphs = '<w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">The transcript features Alberto discussing an app designed for individuals in their 20s and 30s who find themselves spending excessive time on their smartphones. Here are the key points and concepts:</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="1001"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Target Audience</w:t></w:r><w:r><w:t xml:space="preserve">: The app is specifically aimed at young adults who feel overwhelmed by the time they spend on their phones, particularly on social media and entertainment platforms.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="1001"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Purpose of the App</w:t></w:r><w:r><w:t xml:space="preserve">: It aims to transform smartphones from time-wasting devices into productivity tools. The app offers a more reliable way to restrict phone usage compared to existing solutions like Apple Screen Time.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="1001"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Subscription Model</w:t></w:r><w:r><w:t xml:space="preserve">: The app operates on a subscription basis, targeting individuals who are frustrated with how companies capitalize on their attention and wish to regain control over their smartphone use.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1001"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Key Definition</w:t></w:r><w:r><w:t xml:space="preserve">:</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="Compact"/><w:numPr><w:numId w:val="1002"/><w:ilvl w:val="1"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Productivity Tool</w:t></w:r><w:r><w:t xml:space="preserve">: A tool or software that helps users manage their time and tasks more efficiently, turning distractions into productive activities.</w:t></w:r></w:p><w:p><w:pPr><w:numPr><w:numId w:val="1001"/><w:ilvl w:val="0"/></w:numPr></w:pPr><w:r><w:rPr><w:b/></w:rPr><w:t xml:space="preserve">Example Problem</w:t></w:r><w:r><w:t xml:space="preserve">: A user spends several hours daily on social media and entertainment apps, reducing their productivity. By using this app, they can set limits on these activities, converting their phone into a tool that aids in achieving their personal and professional goals.</w:t></w:r></w:p><w:p><w:pPr><w:pStyle w:val="FirstParagraph"/></w:pPr><w:r><w:t xml:space="preserve">Overall, the app provides a solution for those looking to balance their digital consumption and enhance productivity by minimizing distractions from their smartphones.</w:t></w:r></w:p>'

$bodyStartMarker = "<w:body>"
$bodyEndMarker = "<w:sectPr"

$bodyStartIdx = $docXml.IndexOf($bodyStartMarker) + $bodyStartMarker.Length
$bodyEndIdx = $docXml.IndexOf($bodyEndMarker)

$docBefore = $docXml.Substring(0, $bodyStartIdx)
$docAfter = $docXml.Substring($bodyEndIdx)

$docXml = $docBefore + $newBodyParagraphs + $docAfter

# ---------------------------------------------------------------------------
# 2. Add the new abstractNum (991) definition and the num (1001 / 1002)
#    mappings to the numbering part.
# ---------------------------------------------------------------------------

$newAbstractNum = '<w:abstractNum w:abstractNumId="991"><w:nsid w:val="ea454b4c"/><w:multiLevelType w:val="multilevel"/><w:lvl w:ilvl="0"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8226;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="0"/></w:tabs><w:ind w:left="480" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="1"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8211;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="720"/></w:tabs><w:ind w:left="1200" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="2"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8226;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="1440"/></w:tabs><w:ind w:left="1920" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="3"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8211;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="2160"/></w:tabs><w:ind w:left="2640" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="4"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8226;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="2880"/></w:tabs><w:ind w:left="3360" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="5"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8211;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="3600"/></w:tabs><w:ind w:left="4080" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="6"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8226;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="4320"/></w:tabs><w:ind w:left="4800" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="7"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8211;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="5040"/></w:tabs><w:ind w:left="5520" w:hanging="480"/></w:pPr></w:lvl><w:lvl w:ilvl="8"><w:numFmt w:val="bullet"/><w:lvlText w:val="&#8226;"/><w:lvlJc w:val="left"/><w:pPr><w:tabs><w:tab w:val="num" w:pos="5760"/></w:tabs><w:ind w:left="6240" w:hanging="480"/></w:pPr></w:lvl></w:abstractNum>'

$abstractNumEndMarker = "</w:abstractNum>"
$lastAbstractNumEndIdx = $numXml.LastIndexOf($abstractNumEndMarker) + $abstractNumEndMarker.Length

$numXml = $numXml.Substring(0, $lastAbstractNumEndIdx) + $newAbstractNum + $numXml.Substring($lastAbstractNumEndIdx)

$newNums = '<w:num w:numId="1001"><w:abstractNumId w:val="991"/></w:num><w:num w:numId="1002"><w:abstractNumId w:val="991"/></w:num>'

$numberingEndMarker = "</w:numbering>"
$numberingEndIdx = $numXml.IndexOf($numberingEndMarker)

$numXml = $numXml.Substring(0, $numberingEndIdx) + $newNums + $numXml.Substring($numberingEndIdx)

# ---------------------------------------------------------------------------
# 3. Push only the two modified parts back into the document as a minimal
#    OOXML package so every other part stays untouched.
# ---------------------------------------------------------------------------

$pkgWrapped = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' `
    + '<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData>' + $docXml + '</pkg:xmlData></pkg:part>' `
    + '<pkg:part pkg:name="/word/numbering.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.numbering+xml"><pkg:xmlData>' + $numXml + '</pkg:xmlData></pkg:part>' `
    + '</pkg:package>'

$d.Content.InsertXML($pkgWrapped)
